$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "model_1_7_0"
$ws.Cells.Item(2,2).Value = 0.9716593622110361
$ws.Cells.Item(2,3).Value = 0.9848163966797699
$ws.Cells.Item(2,4).Value = 0.9714960963015005
$ws.Cells.Item(2,5).Value = 0.9765106555253509
$ws.Cells.Item(2,6).Value = 3.417889355224486
$ws.Cells.Item(2,7).Value = 0.7288246990781945
$ws.Cells.Item(2,8).Value = 3.707779310339084
$ws.Cells.Item(2,9).Value = 2.130688651051019

$ws.Cells.Item(3,1).Value = "model_1_7_1"
$ws.Cells.Item(3,2).Value = 0.9754610717835247
$ws.Cells.Item(3,3).Value = 0.98074487231038
$ws.Cells.Item(3,4).Value = 0.969521961585375
$ws.Cells.Item(3,5).Value = 0.9740378568764824
$ws.Cells.Item(3,6).Value = 2.959402048897043
$ws.Cells.Item(3,7).Value = 0.9242610168431958
$ws.Cells.Item(3,8).Value = 3.964574166710198
$ws.Cells.Item(3,9).Value = 2.354993080796383

$ws.Cells.Item(4,1).Value = "model_1_7_2"
$ws.Cells.Item(4,2).Value = 0.978283291144702
$ws.Cells.Item(4,3).Value = 0.9752237606972762
$ws.Cells.Item(4,4).Value = 0.9671772568083953
$ws.Cells.Item(4,5).Value = 0.9709088005393038
$ws.Cells.Item(4,6).Value = 2.619041553677957
$ws.Cells.Item(4,7).Value = 1.189278643102976
$ws.Cells.Item(4,8).Value = 4.269572666315571
$ws.Cells.Item(4,9).Value = 2.638825813264546

$ws.Cells.Item(5,1).Value = "model_1_7_24"
$ws.Cells.Item(5,2).Value = 0.9800541257333281
$ws.Cells.Item(5,3).Value = 0.864867734621651
$ws.Cells.Item(5,4).Value = 0.9276005173036169
$ws.Cells.Item(5,5).Value = 0.9132842392344228
$ws.Cells.Item(5,6).Value = 2.405478375057995
$ws.Cells.Item(5,7).Value = 6.486453220159458
$ws.Cells.Item(5,8).Value = 9.417703163059464
$ws.Cells.Item(5,9).Value = 7.865876696979679

$ws.Cells.Item(6,1).Value = "model_1_7_23"
$ws.Cells.Item(6,2).Value = 0.9802592904222188
$ws.Cells.Item(6,3).Value = 0.8669152639502347
$ws.Cells.Item(6,4).Value = 0.9283044223530015
$ws.Cells.Item(6,5).Value = 0.9143329132306084
$ws.Cells.Item(6,6).Value = 2.380735452493956
$ws.Cells.Item(6,7).Value = 6.388170229272141
$ws.Cells.Item(6,8).Value = 9.326139403717651
$ws.Cells.Item(6,9).Value = 7.770752808582694

$ws.Cells.Item(7,1).Value = "model_1_7_3"
$ws.Cells.Item(7,2).Value = 0.9803368263828021
$ws.Cells.Item(7,3).Value = 0.9687486985154555
$ws.Cells.Item(7,4).Value = 0.9646014748244854
$ws.Cells.Item(7,5).Value = 0.9673565322058085
$ws.Cells.Item(7,6).Value = 2.371384592562785
$ws.Cells.Item(7,7).Value = 1.500086634239727
$ws.Cells.Item(7,8).Value = 4.604629620229839
$ws.Cells.Item(7,9).Value = 2.961047569254156

$ws.Cells.Item(8,1).Value = "model_1_7_22"
$ws.Cells.Item(8,2).Value = 0.9804806830852195
$ws.Cells.Item(8,3).Value = 0.8691694926151172
$ws.Cells.Item(8,4).Value = 0.9290799057768044
$ws.Cells.Item(8,5).Value = 0.9154878775671746
$ws.Cells.Item(8,6).Value = 2.354035431420698
$ws.Cells.Item(8,7).Value = 6.279965510426024
$ws.Cells.Item(8,8).Value = 9.225264750733231
$ws.Cells.Item(8,9).Value = 7.665987458193869

$ws.Cells.Item(9,1).Value = "model_1_7_21"
$ws.Cells.Item(9,2).Value = 0.9807185674861527
$ws.Cells.Item(9,3).Value = 0.8716489282469794
$ws.Cells.Item(9,4).Value = 0.929933499552772
$ws.Cells.Item(9,5).Value = 0.9167585325368609
$ws.Cells.Item(9,6).Value = 2.325346501842686
$ws.Cells.Item(9,7).Value = 6.160950683038639
$ws.Cells.Item(9,8).Value = 9.114229526384877
$ws.Cells.Item(9,9).Value = 7.550727957179093

$ws.Cells.Item(10,1).Value = "model_1_7_20"
$ws.Cells.Item(10,2).Value = 0.9809728988564196
$ws.Cells.Item(10,3).Value = 0.8743721170008409
$ws.Cells.Item(10,4).Value = 0.9308722444758812
$ws.Cells.Item(10,5).Value = 0.9181549127754154
$ws.Cells.Item(10,6).Value = 2.294674062866275
$ws.Cells.Item(10,7).Value = 6.030235517329467
$ws.Cells.Item(10,8).Value = 8.992117866157329
$ws.Cells.Item(10,9).Value = 7.424064076454329

$ws.Cells.Item(11,1).Value = "model_1_7_19"
$ws.Cells.Item(11,2).Value = 0.9812431181478899
$ws.Cells.Item(11,3).Value = 0.8773588081325215
$ws.Cells.Item(11,4).Value = 0.9319028726043962
$ws.Cells.Item(11,5).Value = 0.9196871468803179
$ws.Cells.Item(11,6).Value = 2.2620855358624
$ws.Cells.Item(11,7).Value = 5.88687203374936
$ws.Cells.Item(11,8).Value = 8.858054065915027
$ws.Cells.Item(11,9).Value = 7.285077063785983

$ws.Cells.Item(12,1).Value = "model_1_7_18"
$ws.Cells.Item(12,2).Value = 0.9815281707907549
$ws.Cells.Item(12,3).Value = 0.8806295013575519
$ws.Cells.Item(12,4).Value = 0.933032945786778
$ws.Cells.Item(12,5).Value = 0.9213661432553955
$ws.Cells.Item(12,6).Value = 2.227708102263976
$ws.Cells.Item(12,7).Value = 5.729876230102814
$ws.Cells.Item(12,8).Value = 8.71105448266059
$ws.Cells.Item(12,9).Value = 7.132777431695592

$ws.Cells.Item(13,1).Value = "model_1_7_4"
$ws.Cells.Item(13,2).Value = 0.9817902060680403
$ws.Cells.Item(13,3).Value = 0.9616992650088036
$ws.Cells.Item(13,4).Value = 0.9619024732891264
$ws.Cells.Item(13,5).Value = 0.9635602287494021
$ws.Cells.Item(13,6).Value = 2.196106569807434
$ws.Cells.Item(13,7).Value = 1.838464893062646
$ws.Cells.Item(13,8).Value = 4.955714936726491
$ws.Cells.Item(13,9).Value = 3.305405441788265

$ws.Cells.Item(14,1).Value = "model_1_7_17"
$ws.Cells.Item(14,2).Value = 0.9818261270113613
$ws.Cells.Item(14,3).Value = 0.8842049774142495
$ws.Cells.Item(14,4).Value = 0.934270186592942
$ws.Cells.Item(14,5).Value = 0.9232026872530835
$ws.Cells.Item(14,6).Value = 2.191774493348163
$ws.Cells.Item(14,7).Value = 5.558250614883276
$ws.Cells.Item(14,8).Value = 8.550114566797626
$ws.Cells.Item(14,9).Value = 6.966186854540361

$ws.Cells.Item(15,1).Value = "model_1_7_16"
$ws.Cells.Item(15,2).Value = 0.9821339803896
$ws.Cells.Item(15,3).Value = 0.8881043373061259
$ws.Cells.Item(15,4).Value = 0.9356216860246738
$ws.Cells.Item(15,5).Value = 0.9252072389566453
$ws.Cells.Item(15,6).Value = 2.154647284275202
$ws.Cells.Item(15,7).Value = 5.371078325153612
$ws.Cells.Item(15,8).Value = 8.3743119229242
$ws.Cells.Item(15,9).Value = 6.784356511431659

$ws.Cells.Item(16,1).Value = "model_1_7_15"
$ws.Cells.Item(16,2).Value = 0.9824473786204768
$ws.Cells.Item(16,3).Value = 0.8923471621862017
$ws.Cells.Item(16,4).Value = 0.9370951306360364
$ws.Cells.Item(16,5).Value = 0.9273900803298957
$ws.Cells.Item(16,6).Value = 2.116851364323202
$ws.Cells.Item(16,7).Value = 5.167419450429012
$ws.Cells.Item(16,8).Value = 8.182646686375257
$ws.Cells.Item(16,9).Value = 6.586353738470117

$ws.Cells.Item(17,1).Value = "model_1_7_14"
$ws.Cells.Item(17,2).Value = 0.9827603018442608
$ws.Cells.Item(17,3).Value = 0.8969507144418781
$ws.Cells.Item(17,4).Value = 0.9386970811675567
$ws.Cells.Item(17,5).Value = 0.9297608767426826
$ws.Cells.Item(17,6).Value = 2.079112730368037
$ws.Cells.Item(17,7).Value = 4.946445382767233
$ws.Cells.Item(17,8).Value = 7.974265437975578
$ws.Cells.Item(17,9).Value = 6.371301802213261

$ws.Cells.Item(18,1).Value = "model_1_7_5"
$ws.Cells.Item(18,2).Value = 0.9827780740769174
$ws.Cells.Item(18,3).Value = 0.9543630733737788
$ws.Cells.Item(18,4).Value = 0.9591626388857829
$ws.Cells.Item(18,5).Value = 0.9596560485181792
$ws.Cells.Item(18,6).Value = 2.076969393818305
$ws.Cells.Item(18,7).Value = 2.190607763764036
$ws.Cells.Item(18,8).Value = 5.312111780538692
$ws.Cells.Item(18,9).Value = 3.659548679770158

$ws.Cells.Item(19,1).Value = "model_1_7_13"
$ws.Cells.Item(19,2).Value = 0.9830644498226336
$ws.Cells.Item(19,3).Value = 0.9019292724863923
$ws.Cells.Item(19,4).Value = 0.9404340966277315
$ws.Cells.Item(19,5).Value = 0.9323279294874226
$ws.Cells.Item(19,6).Value = 2.042432393622112
$ws.Cells.Item(19,7).Value = 4.707470747293065
$ws.Cells.Item(19,8).Value = 7.748314983851841
$ws.Cells.Item(19,9).Value = 6.138447702952643

$ws.Cells.Item(20,1).Value = "model_1_7_12"
$ws.Cells.Item(20,2).Value = 0.9833485627855617
$ws.Cells.Item(20,3).Value = 0.907292644421922
$ws.Cells.Item(20,4).Value = 0.942311009251482
$ws.Cells.Item(20,5).Value = 0.9350969318102217
$ws.Cells.Item(20,6).Value = 2.008168285703856
$ws.Cells.Item(20,7).Value = 4.450024747518507
$ws.Cells.Item(20,8).Value = 7.504166748323563
$ws.Cells.Item(20,9).Value = 5.887275013553438

$ws.Cells.Item(21,1).Value = "model_1_7_6"
$ws.Cells.Item(21,2).Value = 0.9834073627166362
$ws.Cells.Item(21,3).Value = 0.9469559630897382
$ws.Cells.Item(21,4).Value = 0.9564432814882547
$ws.Cells.Item(21,5).Value = 0.9557457098326577
$ws.Cells.Item(21,6).Value = 2.001076996509724
$ws.Cells.Item(21,7).Value = 2.546154784451291
$ws.Cells.Item(21,8).Value = 5.665844981528398
$ws.Cells.Item(21,9).Value = 4.014250543332102

$ws.Cells.Item(22,1).Value = "model_1_7_11"
$ws.Cells.Item(22,2).Value = 0.9835978729382344
$ws.Cells.Item(22,3).Value = 0.9130443212462059
$ws.Cells.Item(22,4).Value = 0.9443305197724386
$ws.Cells.Item(22,5).Value = 0.9380712330437556
$ws.Cells.Item(22,6).Value = 1.978101407064259
$ws.Cells.Item(22,7).Value = 4.173939812852931
$ws.Cells.Item(22,8).Value = 7.241469420763858
$ws.Cells.Item(22,9).Value = 5.617479920295843

$ws.Cells.Item(23,1).Value = "model_1_7_7"
$ws.Cells.Item(23,2).Value = 0.9837631836973593
$ws.Cells.Item(23,3).Value = 0.9396361975231755
$ws.Cells.Item(23,4).Value = 0.9537897467301758
$ws.Cells.Item(23,5).Value = 0.951904418112702
$ws.Cells.Item(23,6).Value = 1.95816488029573
$ws.Cells.Item(23,7).Value = 2.897509191166134
$ws.Cells.Item(23,8).Value = 6.011015993167353
$ws.Cells.Item(23,9).Value = 4.36268924420428

$ws.Cells.Item(24,1).Value = "model_1_7_10"
$ws.Cells.Item(24,2).Value = 0.9837929881258707
$ws.Cells.Item(24,3).Value = 0.9191798639985141
$ws.Cells.Item(24,4).Value = 0.9464936275046365
$ws.Cells.Item(24,5).Value = 0.9412497785252237
$ws.Cells.Item(24,6).Value = 1.95457045734356
$ws.Cells.Item(24,7).Value = 3.879429016843496
$ws.Cells.Item(24,8).Value = 6.96009301070043
$ws.Cells.Item(24,9).Value = 5.329158090305105

$ws.Cells.Item(25,1).Value = "model_1_7_9"
$ws.Cells.Item(25,2).Value = 0.9839087252070432
$ws.Cells.Item(25,3).Value = 0.9256817332688343
$ws.Cells.Item(25,4).Value = 0.9487971509177602
$ws.Cells.Item(25,5).Value = 0.944625924085317
$ws.Cells.Item(25,6).Value = 1.940612530895684
$ws.Cells.Item(25,7).Value = 3.567334264731977
$ws.Cells.Item(25,8).Value = 6.660451370649861
$ws.Cells.Item(25,9).Value = 5.022912207753948

$ws.Cells.Item(26,1).Value = "model_1_7_8"
$ws.Cells.Item(26,2).Value = 0.98391263573567
$ws.Cells.Item(26,3).Value = 0.932518281531611
$ws.Cells.Item(26,4).Value = 0.9512338685528039
$ws.Cells.Item(26,5).Value = 0.9481855131970014
$ws.Cells.Item(26,6).Value = 1.940140919979011
$ws.Cells.Item(26,7).Value = 3.239174662214366
$ws.Cells.Item(26,8).Value = 6.343483865850553
$ws.Cells.Item(26,9).Value = 4.700026393258097

Write-Host "done"